$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, exactly as scraped/displayed on the
# coinranking.com cryptos list for this run (price in column D, 1h change in E).
$updates = [ordered]@{
    "D2" = "35.422.95"
    "E2" = "  +2.71%  "
    "D3" = "1.836.55"
    "E3" = "  +1.58%  "
    "E4" = "  +0.34%  "
    "D5" = "230.64"
    "E5" = "  +2.48%  "
    "D6" = "0.610"
    "E6" = "  +1.52%  "
    "E7" = "  +0.27%  "
    "D8" = "43.85"
    "E8" = "  +13.80%  "
    "D9" = "0.308"
    "E9" = "  +7.03%  "
    "D10" = "0.0700"
    "E10" = "  +4.48%  "
    "E11" = "  +2.62%  "
    "D12" = "2.102.85"
    "E12" = "  +1.59%  "
    "D13" = "1.829.05"
    "E13" = "  +1.21%  "
    "D14" = "11.26"
    "E14" = "  +1.28%  "
    "D15" = "0.672"
    "E15" = "  +6.60%  "
    "D16" = "4.68"
    "E16" = "  +6.76%  "
    "D17" = "35.390.08"
    "E17" = "  +2.67%  "
    "D18" = "69.91"
    "E18" = "  +2.55%  "
    "E19" = "  +3.88%  "
    "D20" = "244.13"
    "E20" = "  +1.15%  "
    "D21" = "12.03"
    "E21" = "  +7.61%  "
    "D22" = "4.68"
    "E22" = "  +14.13%  "
    "E23" = "  +0.28%  "
    "D24" = "2.19"
    "E24" = "  +0.47%  "
    "D25" = "169.65"
    "E25" = "  -0.53%  "
    "D26" = "7.90"
    "E26" = "  +2.91%  "
    "D27" = "17.69"
    "E27" = "  +1.17%  "
    "E28" = "  -0.43%  "
    "E29" = "  +22.68%  "
    "E30" = "  +0.24%  "
    "D31" = "3.320.02"
    "E31" = "  +36.64%  "
    "D32" = "0.0550"
    "E32" = "  +6.92%  "
    "D33" = "4.08"
    "E33" = "  +6.20%  "
    "D34" = "3.93"
    "E34" = "  +4.05%  "
    "E35" = "  +0.94%  "
    "D36" = "96.04"
    "E36" = "  +16.16%  "
    "E37" = "  +6.85%  "
    "E38" = "  +3.40%  "
    "D39" = "1.344.98"
    "E39" = "  +2.94%  "
    "D40" = "15.49"
    "E40" = "  +10.11%  "
    "E41" = "  +3.99%  "
    "D42" = "2.41"
    "E42" = "  +3.96%  "
    "E43" = "  +6.45%  "
    "E44" = "  +4.06%  "
    "E45" = "  +0.76%  "
    "E46" = "  -0.04%  "
    "D47" = "6.24"
    "E47" = "  +7.51%  "
    "D48" = "0.0520"
    "E48" = "  +0.98%  "
    "D49" = "2.004.72"
    "E49" = "  +1.71%  "
    "E50" = "  +0.29%  "
    "D51" = "102.74"
    "E51" = "  -0.18%  "
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    # Values that look like plain numbers (e.g. "230.64", "0.0700") must be
    # forced to text so Excel keeps the original string (and trailing zeros)
    # instead of silently converting the cell to a number.
    if ($value -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
